$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add four new header cells (M1:P1) for the per-source systematic
# uncertainty breakdown columns, using the same centered/top-aligned
# style used elsewhere for this workbook's header-like cells.
$ws.Range("M1").Value = "syst0_c"
$ws.Range("N1").Value = "syst1_c"
$ws.Range("O1").Value = "syst2_c"
$ws.Range("P1").Value = "syst3_c"

# Rename the "syst_u" header (column H) to "syst_tot" - the old "syst_u"
# shared string will be dropped automatically once it becomes unused.
$ws.Range("H1").Value = "syst_tot"

$headerRange = $ws.Range("M1:P1")
$headerRange.Font.Size = 11
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160


# Fill in the new systematic-uncertainty component values for each data row.
$ws.Range("M2").Value = 0.0018
$ws.Range("N2").Value = 0.0006
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.0045

$ws.Range("M3").Value = 0.0025
$ws.Range("N3").Value = 0.0006
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.0044

$ws.Range("M4").Value = 0.0027
$ws.Range("N4").Value = 0.0003
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0044

$ws.Range("M5").Value = 0.0025
$ws.Range("N5").Value = 0.0003
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0.0044

$ws.Range("M6").Value = 0.0019
$ws.Range("N6").Value = 0.0006
$ws.Range("O6").Value = 0.0001
$ws.Range("P6").Value = 0.0044

$ws.Range("M7").Value = 0.0024
$ws.Range("N7").Value = 0.001
$ws.Range("O7").Value = 0.0001
$ws.Range("P7").Value = 0.0049

$ws.Range("M8").Value = 0.0026
$ws.Range("N8").Value = 0.0008
$ws.Range("O8").Value = 0.0001
$ws.Range("P8").Value = 0.0054

$ws.Range("M9").Value = 0.0031
$ws.Range("N9").Value = 0.0008
$ws.Range("O9").Value = 0.0001
$ws.Range("P9").Value = 0.0092

$ws.Range("M10").Value = 0.002
$ws.Range("N10").Value = 0.0016
$ws.Range("O10").Value = 0.0002
$ws.Range("P10").Value = 0.0087

$ws.Range("M11").Value = 0.002
$ws.Range("N11").Value = 0.0026
$ws.Range("O11").Value = 0.0003
$ws.Range("P11").Value = 0.01

$ws.Range("M12").Value = 0.0029
$ws.Range("N12").Value = 0.0024
$ws.Range("O12").Value = 0.0003
$ws.Range("P12").Value = 0.0125

# Mirror the author's final cell selection.
$ws.Range("L18").Select()
